$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3067331670822943
$ws.Range("B3").Value = 0.6034912718204489
$ws.Range("B4").Value = 0.6209476309226932
$ws.Range("B5").Value = 0.2992518703241895
